$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text  = "415×2="
$t.Cell(1, 2).Range.Text  = "489×4="
$t.Cell(1, 3).Range.Text  = "520×3="
$t.Cell(1, 4).Range.Text  = "882×3="
$t.Cell(1, 5).Range.Text  = "251×6="

$t.Cell(5, 1).Range.Text  = "613×2="
$t.Cell(5, 2).Range.Text  = "869×9="
$t.Cell(5, 3).Range.Text  = "811×7="
$t.Cell(5, 4).Range.Text  = "135×3="
$t.Cell(5, 5).Range.Text  = "527×6="

$t.Cell(10, 1).Range.Text = "772×8="
$t.Cell(10, 2).Range.Text = "448×8="
$t.Cell(10, 3).Range.Text = "146×5="
$t.Cell(10, 4).Range.Text = "909×2="
$t.Cell(10, 5).Range.Text = "624×4="

$t.Cell(15, 1).Range.Text = "434×5="
$t.Cell(15, 2).Range.Text = "633×2="
$t.Cell(15, 3).Range.Text = "253×4="
$t.Cell(15, 4).Range.Text = "955×8="
$t.Cell(15, 5).Range.Text = "176×3="

$t.Cell(20, 1).Range.Text = "342×2="
$t.Cell(20, 2).Range.Text = "173×2="
$t.Cell(20, 3).Range.Text = "675×7="
$t.Cell(20, 4).Range.Text = "582×6="
$t.Cell(20, 5).Range.Text = "962×5="
